$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "RM 232" row (row 26) and the "SC 92" row (originally row 28).
# Delete from the bottom up so row indices of earlier rows remain valid.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# After the deletions, rows have shifted up by two starting at row 28 (old row
# numbering), and by one between rows 27..27. The surviving rows (by former
# row number) now sit at:
#   27 (SC 5)   -> 26
#   29 (SC 101) -> 27
#   30 (SC 105) -> 28
#   31 (SC 119) -> 29
#   32 (SC 120) -> 30
#   33 (SC 132) -> 31
#   34 (SC 193) -> 32
#   35 (SC 232) -> 33
#
# Update column E ("D" header) values that moved between rows as part of the
# re-sampling of which points are treated as missing.
$ws.Range("E27").Value = -10
$ws.Range("E28").Value = ""
$ws.Range("E29").Value = ""
$ws.Range("E30").Value = -5.7
$ws.Range("E32").Value = ""
